$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header row had two placeholder/"unnamed" column labels left over from a
# pandas multi-index export ("unnamed: 1_level_1" and "unnamed: 5_level_1").
# Correct them to read "total", matching the other "total" column header.
$ws.Range("B2").Value = "total"
$ws.Range("F2").Value = "total"
